$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28 (anchor 27772)
$ws.Range("H28").Value = 939.4666999999999
$ws.Range("I28").Value = 354.2857
$ws.Range("K28").Value = 354.2857
$ws.Range("M28").Value = 130.7143
# Row 132 (anchor 44049)
$ws.Range("H132").Value = 4770.4443
$ws.Range("I132").Value = 4770.4443
$ws.Range("K132").Value = 14311.3329
$ws.Range("M132").Value = -11781.3329
# Row 135 (anchor 44047)
$ws.Range("H135").Value = 16134207
$ws.Range("I135").Value = 616.2727
$ws.Range("K135").Value = 5546.454299999999
$ws.Range("M135").Value = -3011.454299999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (anchor 44147)
$ws.Range("H32").Value = 7316.259
$ws.Range("I32").Value = 4833.245
$ws.Range("K32").Value = 4833.245
$ws.Range("M32").Value = -4546.245
# Row 44 (anchor 3861)
$ws.Range("H44").Value = 31024.5
$ws.Range("J44").Value = 31024.5
$ws.Range("L44").Value = 31024.5
$ws.Range("N44").Value = -32000.5
# Row 55 (anchor 2830)
$ws.Range("H55").Value = 25439.75
$ws.Range("J55").Value = 25439.75
$ws.Range("L55").Value = 25439.75
$ws.Range("N55").Value = -26069.75
# Row 98 (anchor 18371)
$ws.Range("H98").Value = 29998.5
$ws.Range("J98").Value = 29998.5
$ws.Range("L98").Value = 29998.5
$ws.Range("N98").Value = -35988.5
# Row 122 (anchor 36168)
$ws.Range("H122").Value = 2218.2354
$ws.Range("I122").Value = 1715.0714
$ws.Range("K122").Value = 5145.2142
$ws.Range("M122").Value = -2695.2142
# Row 132 (anchor 43997)
$ws.Range("H132").Value = 17615.781
$ws.Range("I132").Value = 1715.1154
$ws.Range("J132").Value = 86518.664
$ws.Range("K132").Value = 5145.3462
$ws.Range("L132").Value = 259555.992
$ws.Range("M132").Value = -2615.3462
$ws.Range("N132").Value = -264615.992

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 11 (anchor 2481)
$ws.Range("H11").Value = 1999
$ws.Range("I11").Value = 1999
$ws.Range("K11").Value = 1999
$ws.Range("M11").Value = -1859
# Row 134 (anchor 43998)
$ws.Range("H134").Value = 2832.587
$ws.Range("I134").Value = 2969.244
$ws.Range("K134").Value = 8907.732
$ws.Range("M134").Value = -6372.732

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58 (anchor 44021)
$ws.Range("H58").Value = 21727.44
$ws.Range("I58").Value = 1674.5714
$ws.Range("J58").Value = 47249.273
$ws.Range("K58").Value = 1674.5714
$ws.Range("L58").Value = 47249.273
$ws.Range("M58").Value = -1471.5714
$ws.Range("N58").Value = -47655.273
# Row 95 (anchor 18192)
$ws.Range("H95").Value = 33333
$ws.Range("J95").Value = 33333
$ws.Range("L95").Value = 33333
$ws.Range("N95").Value = -38825
# Row 132 (anchor 44019)
$ws.Range("H132").Value = 3187.625
$ws.Range("I132").Value = 2415.1333
$ws.Range("J132").Value = 4475.1113
$ws.Range("K132").Value = 7245.3999
$ws.Range("L132").Value = 13425.3339
$ws.Range("M132").Value = -4715.3999
$ws.Range("N132").Value = -18485.3339
# Row 136 (anchor 44021)
$ws.Range("H136").Value = 21727.44
$ws.Range("I136").Value = 1674.5714
$ws.Range("J136").Value = 47249.273
$ws.Range("K136").Value = 5023.7142
$ws.Range("L136").Value = 141747.819
$ws.Range("M136").Value = -2473.7142
$ws.Range("N136").Value = -146847.819

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7 (anchor 4728)
$ws.Range("H7").Value = 186.75
$ws.Range("I7").Value = 131
$ws.Range("J7").Value = 205.33333
$ws.Range("K7").Value = 393
$ws.Range("L7").Value = 615.99999
$ws.Range("M7").Value = -281
$ws.Range("N7").Value = -839.99999
# Row 75 (anchor 12863)
$ws.Range("H75").Value = 738.2
$ws.Range("J75").Value = 668.6667
$ws.Range("L75").Value = 2006.0001
$ws.Range("N75").Value = -4002.0001
# Row 78 (anchor 12863)
$ws.Range("H78").Value = 738.2
$ws.Range("J78").Value = 668.6667
$ws.Range("L78").Value = 6018.0003
$ws.Range("N78").Value = -16002.0003

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 39 (anchor 18264)
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
# Row 95 (anchor 18235)
$ws.Range("H95").Value = 21562.666
$ws.Range("J95").Value = 21562.666
$ws.Range("L95").Value = 21562.666
$ws.Range("N95").Value = -27054.666
# Row 102 (anchor 36169)
$ws.Range("H102").Value = 26318922
$ws.Range("I102").Value = 38464880
$ws.Range("J102").Value = 2677.3333
$ws.Range("K102").Value = 38464880
$ws.Range("L102").Value = 2677.3333
$ws.Range("M102").Value = -38463258
$ws.Range("N102").Value = -5921.3333
# Row 107 (anchor 27802)
$ws.Range("H107").Value = 5494754
$ws.Range("J107").Value = 19230890
$ws.Range("L107").Value = 19230890
$ws.Range("N107").Value = -19234730
# Row 113 (anchor 27710)
$ws.Range("H113").Value = 3414.1428
$ws.Range("I113").Value = 2799.6667
$ws.Range("K113").Value = 2799.6667
$ws.Range("M113").Value = -629.6667000000002
# Row 132 (anchor 44008)
$ws.Range("H132").Value = 28810.158
$ws.Range("I132").Value = 2095.7273
$ws.Range("J132").Value = 65542.5
$ws.Range("K132").Value = 6287.1819
$ws.Range("L132").Value = 196627.5
$ws.Range("M132").Value = -3757.1819
$ws.Range("N132").Value = -201687.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40 (anchor 36248)
$ws.Range("H40").Value = 4068.875
$ws.Range("I40").Value = 2984.3333
$ws.Range("K40").Value = 2984.3333
$ws.Range("M40").Value = -2848.3333
# Row 46 (anchor 5282)
$ws.Range("H46").Value = 1028.2069
$ws.Range("I46").Value = 1000.76
$ws.Range("J46").Value = 1199.75
$ws.Range("K46").Value = 1000.76
$ws.Range("L46").Value = 1199.75
$ws.Range("M46").Value = -812.76
$ws.Range("N46").Value = -1575.75
# Row 61 (anchor 27740)
$ws.Range("H61").Value = 3951.8125
$ws.Range("J61").Value = 9347.4
$ws.Range("L61").Value = 9347.4
$ws.Range("N61").Value = -9751.4
# Row 98 (anchor 18379)
$ws.Range("H98").Value = 23000
$ws.Range("J98").Value = 23000
$ws.Range("L98").Value = 23000
$ws.Range("N98").Value = -28990
# Row 113 (anchor 27740)
$ws.Range("H113").Value = 3951.8125
$ws.Range("J113").Value = 9347.4
$ws.Range("L113").Value = 9347.4
$ws.Range("N113").Value = -13687.4
# Row 122 (anchor 36247)
$ws.Range("H122").Value = 787182.4399999999
$ws.Range("I122").Value = 1785063
$ws.Range("J122").Value = 3133.4285
$ws.Range("K122").Value = 5355189
$ws.Range("L122").Value = 9400.2855
$ws.Range("M122").Value = -5352739
$ws.Range("N122").Value = -14300.2855
# Row 132 (anchor 44058)
$ws.Range("H132").Value = 2052.9565
$ws.Range("I132").Value = 1577.1428
$ws.Range("J132").Value = 2261.125
$ws.Range("K132").Value = 4731.428400000001
$ws.Range("L132").Value = 6783.375
$ws.Range("M132").Value = -2201.428400000001
$ws.Range("N132").Value = -11843.375

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132 (anchor 44029)
$ws.Range("H132").Value = 1023.381
$ws.Range("I132").Value = 602.4706
$ws.Range("K132").Value = 1807.4118
$ws.Range("M132").Value = 722.5882000000001
# Row 136 (anchor 44031)
$ws.Range("H136").Value = 23258238
$ws.Range("I136").Value = 34484030
$ws.Range("K136").Value = 103452090
$ws.Range("M136").Value = -103449540